# Update cryptocurrency price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D2").Value = "264.48"
$ws.Range("D3").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D3").Value = "22.82"
$ws.Range("D4").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D4").Value = "6.233"
$ws.Range("D5").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D5").Value = "0.06126"
$ws.Range("D6").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D6").Value = "3.558"
$ws.Range("D7").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D7").Value = "6.732"
$ws.Range("D8").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D8").Value = "1.375"
$ws.Range("D9").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D9").Value = "0.8142"
$ws.Range("D10").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D10").Value = "0.1594"
$ws.Range("D11").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D11").Value = "0.08214"
$ws.Range("D12").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D12").Value = "0.03391"
$ws.Range("D13").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D13").Value = "0.03174"
$ws.Range("D14").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D14").Value = "0.09247"
$ws.Range("D15").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D15").Value = "3.921"
$ws.Range("D16").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D16").Value = "0.001691"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("D17").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D17").Value = "0.04864"
$ws.Range("D18").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D18").Value = "0.0006272"
$ws.Range("D19").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D19").Value = "0.006238"
$ws.Range("D20").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D20").Value = "0.001107"
$ws.Range("D21").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D21").Value = "0.003206"
$ws.Range("D22").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D22").Value = "0.0001507"
$ws.Range("D23").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D23").Value = "3.691"
$ws.Range("D24").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D24").Value = "2.260"
$ws.Range("D25").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D25").Value = "0.3387"
$ws.Range("D40").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D40").Value = "0.04583"
$ws.Range("D41").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D41").Value = "0.1128"
$ws.Range("D42").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D42").Value = "0.003144"
$ws.Range("D43").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D43").Value = "0.003451"
$ws.Range("D44").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D44").Value = "0.01080"
$ws.Range("D45").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D45").Value = "0.00006178"
$ws.Range("D47").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D47").Value = "0.7525"
$ws.Range("D48").NumberFormat = "@"  # keep numeric-looking text as Text, like the source cell
$ws.Range("D48").Value = "0.2488"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
